$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2568.644
$ws.Range("I17").Value = 850
$ws.Range("J17").Value = 2628.9473
$ws.Range("K17").Value = 2550
$ws.Range("L17").Value = 7886.841899999999
$ws.Range("M17").Value = -2382
$ws.Range("N17").Value = -8222.841899999999
$ws.Range("H19").Value = 4412.5713
$ws.Range("I19").Value = 5797.2
$ws.Range("K19").Value = 5797.2
$ws.Range("M19").Value = -5622.2
$ws.Range("H40").Value = 4227.0303
$ws.Range("I40").Value = 3907.6155
$ws.Range("J40").Value = 4434.65
$ws.Range("K40").Value = 3907.6155
$ws.Range("L40").Value = 4434.65
$ws.Range("M40").Value = -3732.6155
$ws.Range("N40").Value = -4784.65
$ws.Range("H70").Value = 6412.5
$ws.Range("J70").Value = 7688.75
$ws.Range("L70").Value = 23066.25
$ws.Range("N70").Value = -23606.25
$ws.Range("H73").Value = 6412.5
$ws.Range("J73").Value = 7688.75
$ws.Range("L73").Value = 23066.25
$ws.Range("N73").Value = -24938.25
$ws.Range("H86").Value = 5329.6665
$ws.Range("J86").Value = 5494.5
$ws.Range("L86").Value = 5494.5
$ws.Range("N86").Value = -7740.5
$ws.Range("H89").Value = 5329.6665
$ws.Range("J89").Value = 5494.5
$ws.Range("L89").Value = 27472.5
$ws.Range("N89").Value = -38704.5
$ws.Range("H101").Value = 1280.0714
$ws.Range("I101").Value = 660.2222
$ws.Range("J101").Value = 2395.8
$ws.Range("K101").Value = 1980.6666
$ws.Range("L101").Value = 7187.400000000001
$ws.Range("M101").Value = -358.6666
$ws.Range("N101").Value = -10431.4
$ws.Range("H121").Value = 2949.2856
$ws.Range("J121").Value = 2949.2856
$ws.Range("L121").Value = 8847.856800000001
$ws.Range("N121").Value = -12341.8568
$ws.Range("H129").Value = 1896.76
$ws.Range("I129").Value = 933.1539
$ws.Range("K129").Value = 2799.4617
$ws.Range("M129").Value = 2200.5383
$ws.Range("H132").Value = 940.8246
$ws.Range("I132").Value = 793
$ws.Range("K132").Value = 2379
$ws.Range("M132").Value = 151
$ws.Range("H138").Value = 3244.34
$ws.Range("I138").Value = 2284.2307
$ws.Range("K138").Value = 6852.6921
$ws.Range("M138").Value = -1712.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2907.4119
$ws.Range("I32").Value = 2107.3547
$ws.Range("K32").Value = 2107.3547
$ws.Range("M32").Value = -1820.3547
$ws.Range("H75").Value = 70989.5
$ws.Range("J75").Value = 70989.5
$ws.Range("L75").Value = 70989.5
$ws.Range("N75").Value = -72737.5
$ws.Range("H78").Value = 70989.5
$ws.Range("J78").Value = 70989.5
$ws.Range("L78").Value = 212968.5
$ws.Range("N78").Value = -221704.5
$ws.Range("H93").Value = 30203.5
$ws.Range("I93").Value = 25407
$ws.Range("K93").Value = 25407
$ws.Range("M93").Value = -22911
$ws.Range("H122").Value = 2062.2144
$ws.Range("I122").Value = 1198
$ws.Range("J122").Value = 2710.375
$ws.Range("K122").Value = 3594
$ws.Range("L122").Value = 8131.125
$ws.Range("M122").Value = -1144
$ws.Range("N122").Value = -13031.125
$ws.Range("H132").Value = 2086.3818
$ws.Range("I132").Value = 1655.415
$ws.Range("K132").Value = 4966.245
$ws.Range("M132").Value = -2436.245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1889.4286
$ws.Range("I99").Value = 1737.8462
$ws.Range("J99").Value = 2135.75
$ws.Range("K99").Value = 1737.8462
$ws.Range("L99").Value = 2135.75
$ws.Range("M99").Value = -239.8462
$ws.Range("N99").Value = -5131.75
$ws.Range("H126").Value = 1889.4286
$ws.Range("I126").Value = 1737.8462
$ws.Range("J126").Value = 2135.75
$ws.Range("K126").Value = 5213.5386
$ws.Range("L126").Value = 6407.25
$ws.Range("M126").Value = -2743.5386
$ws.Range("N126").Value = -11347.25
$ws.Range("H132").Value = 14667.833
$ws.Range("I132").Value = 16748.25
$ws.Range("J132").Value = 10507
$ws.Range("K132").Value = 50244.75
$ws.Range("L132").Value = 31521
$ws.Range("M132").Value = -47714.75
$ws.Range("N132").Value = -36581

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 29042180
$ws.Range("I131").Value = 20835724
$ws.Range("K131").Value = 62507172
$ws.Range("M131").Value = -62502132
$ws.Range("H134").Value = 893
$ws.Range("I134").Value = 893
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2679
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 2391
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 14323.667
$ws.Range("I29").Value = 14273.429
$ws.Range("J29").Value = 14499.5
$ws.Range("K29").Value = 14273.429
$ws.Range("L29").Value = 14499.5
$ws.Range("M29").Value = -13983.429
$ws.Range("N29").Value = -15079.5
$ws.Range("H107").Value = 1138.2727
$ws.Range("I107").Value = 445.25
$ws.Range("J107").Value = 1969.9
$ws.Range("K107").Value = 445.25
$ws.Range("L107").Value = 1969.9
$ws.Range("M107").Value = 1474.75
$ws.Range("N107").Value = -5809.9
$ws.Range("H113").Value = 3749.5
$ws.Range("I113").Value = 3054.3845
$ws.Range("J113").Value = 5040.4287
$ws.Range("K113").Value = 3054.3845
$ws.Range("L113").Value = 5040.4287
$ws.Range("M113").Value = -884.3845000000001
$ws.Range("N113").Value = -9380.4287
$ws.Range("H118").Value = 25894.25
$ws.Range("J118").Value = 25894.25
$ws.Range("L118").Value = 25894.25
$ws.Range("N118").Value = -29208.25
$ws.Range("H122").Value = 4280.0454
$ws.Range("I122").Value = 3757.65
$ws.Range("K122").Value = 11272.95
$ws.Range("M122").Value = -8822.950000000001
$ws.Range("H127").Value = 52378.43
$ws.Range("J127").Value = 52378.43
$ws.Range("L127").Value = 52378.43
$ws.Range("N127").Value = -62298.43
$ws.Range("H135").Value = 69998.14
$ws.Range("J135").Value = 69998.14
$ws.Range("L135").Value = 69998.14
$ws.Range("N135").Value = -80138.14

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9259.235000000001
$ws.Range("I7").Value = 4549.25
$ws.Range("J7").Value = 13445.889
$ws.Range("K7").Value = 4549.25
$ws.Range("L7").Value = 13445.889
$ws.Range("M7").Value = -4437.25
$ws.Range("N7").Value = -13669.889
$ws.Range("H46").Value = 2579.8635
$ws.Range("J46").Value = 2787.5264
$ws.Range("L46").Value = 2787.5264
$ws.Range("N46").Value = -3163.5264
$ws.Range("H126").Value = 9259.235000000001
$ws.Range("I126").Value = 4549.25
$ws.Range("J126").Value = 13445.889
$ws.Range("K126").Value = 13647.75
$ws.Range("L126").Value = 40337.667
$ws.Range("M126").Value = -11177.75
$ws.Range("N126").Value = -45277.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 10000
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10586
$ws.Range("H32").Value = 28029
$ws.Range("J32").Value = 28029
$ws.Range("L32").Value = 28029
$ws.Range("N32").Value = -28663
$ws.Range("H81").Value = 5811.5
$ws.Range("I81").Value = 2593
$ws.Range("K81").Value = 5186
$ws.Range("M81").Value = -4125
$ws.Range("H84").Value = 5811.5
$ws.Range("I84").Value = 2593
$ws.Range("K84").Value = 25930
$ws.Range("M84").Value = -20626
$ws.Range("H93").Value = 49996
$ws.Range("J93").Value = 49996
$ws.Range("L93").Value = 49996
$ws.Range("N93").Value = -54988
$ws.Range("H113").Value = 982
$ws.Range("I113").Value = 1046.1666
$ws.Range("J113").Value = 896.44446
$ws.Range("K113").Value = 3138.4998
$ws.Range("L113").Value = 2689.33338
$ws.Range("M113").Value = -968.4998000000001
$ws.Range("N113").Value = -7029.33338
$ws.Range("H122").Value = 7398.5557
$ws.Range("I122").Value = 2837.524
$ws.Range("K122").Value = 8512.572
$ws.Range("M122").Value = -6062.572
$ws.Range("H126").Value = 2424.7856
$ws.Range("I126").Value = 1285.4286
$ws.Range("J126").Value = 3564.1428
$ws.Range("K126").Value = 3856.2858
$ws.Range("L126").Value = 10692.4284
$ws.Range("M126").Value = -1386.2858
$ws.Range("N126").Value = -15632.4284
$ws.Range("H136").Value = 3914.2273
$ws.Range("I136").Value = 2084.6316
$ws.Range("K136").Value = 6253.8948
$ws.Range("M136").Value = -3703.8948

